{"js": "// Replace the 25 \"division problem\" answer strings in the worksheet table.\n// Each old string is unique in the document, so a direct search + full\n// replace of the matched range is unambiguous and preserves the existing\n// run formatting (font/size) because insertText(\"Replace\") only swaps the\n// text content of the matched range, not the surrounding run properties.\nconst replacements = [\n  [\"676\u00f75=135, 1\", \"819\u00f78=102, 3\"],\n  [\"110\u00f79=12, 2\", \"226\u00f77=32, 2\"],\n  [\"319\u00f77=45, 4\", \"192\u00f76=32, 0\"],\n  [\"102\u00f77=14, 4\", \"473\u00f79=52, 5\"],\n  [\"223\u00f78=27, 7\", \"611\u00f74=152, 3\"],\n  [\"893\u00f76=148, 5\", \"654\u00f72=327, 0\"],\n  [\"786\u00f78=98, 2\", \"415\u00f75=83, 0\"],\n  [\"649\u00f72=324, 1\", \"730\u00f77=104, 2\"],\n  [\"414\u00f72=207, 0\", \"630\u00f72=315, 0\"],\n  [\"490\u00f78=61, 2\", \"536\u00f77=76, 4\"],\n  [\"951\u00f76=158, 3\", \"916\u00f76=152, 4\"],\n  [\"778\u00f78=97, 2\", \"288\u00f74=72, 0\"],\n  [\"584\u00f76=97, 2\", \"734\u00f78=91, 6\"],\n  [\"547\u00f78=68, 3\", \"145\u00f76=24, 1\"],\n  [\"314\u00f79=34, 8\", \"830\u00f78=103, 6\"],\n  [\"977\u00f75=195, 2\", \"906\u00f78=113, 2\"],\n  [\"713\u00f76=118, 5\", \"548\u00f79=60, 8\"],\n  [\"155\u00f79=17, 2\", \"373\u00f73=124, 1\"],\n  [\"741\u00f77=105, 6\", \"490\u00f75=98, 0\"],\n  [\"378\u00f74=94, 2\", \"329\u00f74=82, 1\"],\n  [\"717\u00f75=143, 2\", \"510\u00f72=255, 0\"],\n  [\"342\u00f72=171, 0\", \"642\u00f79=71, 3\"],\n  [\"625\u00f72=312, 1\", \"647\u00f74=161, 3\"],\n  [\"586\u00f78=73, 2\", \"742\u00f79=82, 4\"],\n  [\"251\u00f75=50, 1\", \"309\u00f75=61, 4\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"division problem\" answer strings in the worksheet table.\n# Each old string is unique in the document, so Find/Replace (ReplaceAll)\n# on the exact text is unambiguous and leaves the surrounding run\n# formatting (font/size) untouched.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"676\u00f75=135, 1\", \"819\u00f78=102, 3\"),\n  @(\"110\u00f79=12, 2\", \"226\u00f77=32, 2\"),\n  @(\"319\u00f77=45, 4\", \"192\u00f76=32, 0\"),\n  @(\"102\u00f77=14, 4\", \"473\u00f79=52, 5\"),\n  @(\"223\u00f78=27, 7\", \"611\u00f74=152, 3\"),\n  @(\"893\u00f76=148, 5\", \"654\u00f72=327, 0\"),\n  @(\"786\u00f78=98, 2\", \"415\u00f75=83, 0\"),\n  @(\"649\u00f72=324, 1\", \"730\u00f77=104, 2\"),\n  @(\"414\u00f72=207, 0\", \"630\u00f72=315, 0\"),\n  @(\"490\u00f78=61, 2\", \"536\u00f77=76, 4\"),\n  @(\"951\u00f76=158, 3\", \"916\u00f76=152, 4\"),\n  @(\"778\u00f78=97, 2\", \"288\u00f74=72, 0\"),\n  @(\"584\u00f76=97, 2\", \"734\u00f78=91, 6\"),\n  @(\"547\u00f78=68, 3\", \"145\u00f76=24, 1\"),\n  @(\"314\u00f79=34, 8\", \"830\u00f78=103, 6\"),\n  @(\"977\u00f75=195, 2\", \"906\u00f78=113, 2\"),\n  @(\"713\u00f76=118, 5\", \"548\u00f79=60, 8\"),\n  @(\"155\u00f79=17, 2\", \"373\u00f73=124, 1\"),\n  @(\"741\u00f77=105, 6\", \"490\u00f75=98, 0\"),\n  @(\"378\u00f74=94, 2\", \"329\u00f74=82, 1\"),\n  @(\"717\u00f75=143, 2\", \"510\u00f72=255, 0\"),\n  @(\"342\u00f72=171, 0\", \"642\u00f79=71, 3\"),\n  @(\"625\u00f72=312, 1\", \"647\u00f74=161, 3\"),\n  @(\"586\u00f78=73, 2\", \"742\u00f79=82, 4\"),\n  @(\"251\u00f75=50, 1\", \"309\u00f75=61, 4\")\n)\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Execute($old, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new, $wdReplaceAll) | Out-Null\n}\n"}
